$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking cell updates ---
$ws.Range("D2").Value = "43.125.42"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.572.84"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("D13").Value = "2.965.22"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("E14").Value = "  -3.45%  "
$ws.Range("D15").Value = "2.554.18"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "43.132.13"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").Value = "  +2.99%  "
$ws.Range("E20").Value = "  -3.37%  "
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("E25").Value = "  +2.55%  "
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("E31").Value = "  -3.81%  "
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("E39").Value = "  +6.33%  "
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("E41").Value = "  -5.26%  "
$ws.Range("E42").Value = "  +3.89%  "
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("D46").Value = "1.998.03"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("D48").Value = "2.817.91"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("E49").Value = "  -3.15%  "
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("E51").Value = "  +2.23%  "

# --- Numeric-looking cell updates: force Text storage to preserve original inlineStr formatting ---
# (Excel would otherwise auto-convert these number-like strings into actual numbers)
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "314.94"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "96.74"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.540"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "35.51"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.843"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.84"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.59"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "69.37"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "253.50"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.96"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.08"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "26.95"
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "40.22"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "10.29"
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "5.84"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "154.65"
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.38"
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0807"
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.12"
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.70"
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "18.96"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.43"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "22.39"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.96"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "8.89"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "82.98"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "74.80"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.193"
$c.Style = "Normal"
